$d = $word.ActiveDocument

$d.Content.Find.Execute(", CPF nº ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", RG {{rg}}, CPF nº ", 2)
